$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1268.7778
$ws.Range("I2").Value = 802.5
$ws.Range("J2").Value = 4999
$ws.Range("K2").Value = 802.5
$ws.Range("L2").Value = 4999
$ws.Range("M2").Value = -689.5
$ws.Range("N2").Value = -5225

$ws.Range("H38").Value = 35.2
$ws.Range("I38").Value = 35.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 105.6
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = 266.4

$ws.Range("H70").Value = 1600
$ws.Range("I70").Value = 1766.6666
$ws.Range("K70").Value = 5299.9998
$ws.Range("M70").Value = -5029.9998

$ws.Range("H73").Value = 1600
$ws.Range("I73").Value = 1766.6666
$ws.Range("K73").Value = 5299.9998
$ws.Range("M73").Value = -4363.9998

$ws.Range("H86").Value = 8516.666999999999
$ws.Range("I86").Value = 8025
$ws.Range("K86").Value = 8025
$ws.Range("M86").Value = -6902

$ws.Range("H89").Value = 8516.666999999999
$ws.Range("I89").Value = 8025
$ws.Range("K89").Value = 40125
$ws.Range("M89").Value = -34509

$ws.Range("H132").Value = 2770
$ws.Range("I132").Value = 2806
$ws.Range("K132").Value = 8418
$ws.Range("M132").Value = -5888

$ws.Range("H137").Value = 2249.2
$ws.Range("I137").Value = 1938.8
$ws.Range("K137").Value = 5816.4
$ws.Range("M137").Value = -3266.4

$ws.Range("H141").Value = 1825.8572
$ws.Range("I141").Value = 1825.8572
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5477.571599999999
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -297.5715999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4604.8
$ws.Range("I2").Value = 1728.3334
$ws.Range("K2").Value = 1728.3334
$ws.Range("M2").Value = -1615.3334

$ws.Range("H45").Value = 3078.6924
$ws.Range("I45").Value = 1364.6
$ws.Range("J45").Value = 4150
$ws.Range("K45").Value = 1364.6
$ws.Range("L45").Value = 4150
$ws.Range("M45").Value = -987.5999999999999
$ws.Range("N45").Value = -4904

$ws.Range("H74").Value = 1150.1765
$ws.Range("I74").Value = 1150.1765
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1150.1765
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -276.1765

$ws.Range("H77").Value = 1150.1765
$ws.Range("I77").Value = 1150.1765
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5750.8825
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -1382.8825

$ws.Range("H116").Value = 4604.8
$ws.Range("I116").Value = 1728.3334
$ws.Range("K116").Value = 1728.3334
$ws.Range("M116").Value = 565.6666

$ws.Range("H122").Value = 2644.2856
$ws.Range("I122").Value = 2370.2
$ws.Range("K122").Value = 7110.599999999999
$ws.Range("M122").Value = -4660.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4604.8
$ws.Range("I3").Value = 1728.3334
$ws.Range("K3").Value = 1728.3334
$ws.Range("M3").Value = -1614.3334

$ws.Range("H56").Value = 13749.875
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9261

$ws.Range("H64").Value = 966.3333
$ws.Range("I64").Value = 903.3333
$ws.Range("K64").Value = 903.3333
$ws.Range("M64").Value = -678.3333

$ws.Range("H67").Value = 966.3333
$ws.Range("I67").Value = 903.3333
$ws.Range("K67").Value = 903.3333
$ws.Range("M67").Value = -123.3333

$ws.Range("H99").Value = 5375.5
$ws.Range("I99").Value = 4051
$ws.Range("K99").Value = 4051
$ws.Range("M99").Value = -2553

$ws.Range("H105").Value = 4023
$ws.Range("I105").Value = 3709
$ws.Range("J105").Value = 5750
$ws.Range("K105").Value = 3709
$ws.Range("L105").Value = 5750
$ws.Range("M105").Value = -1962
$ws.Range("N105").Value = -9244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2257
$ws.Range("I31").Value = 2142
$ws.Range("K31").Value = 2142
$ws.Range("M31").Value = -1847

$ws.Range("H34").Value = 2257
$ws.Range("I34").Value = 2142
$ws.Range("K34").Value = 2142
$ws.Range("M34").Value = -1940

$ws.Range("H99").Value = 8041.1665
$ws.Range("I99").Value = 8054.1816
$ws.Range("J99").Value = 7898
$ws.Range("K99").Value = 8054.1816
$ws.Range("L99").Value = 7898
$ws.Range("M99").Value = -6556.1816
$ws.Range("N99").Value = -10894

$ws.Range("H126").Value = 8041.1665
$ws.Range("I126").Value = 8054.1816
$ws.Range("J126").Value = 7898
$ws.Range("K126").Value = 24162.5448
$ws.Range("L126").Value = 23694
$ws.Range("M126").Value = -21692.5448
$ws.Range("N126").Value = -28634

$ws.Range("H132").Value = 3905.7693
$ws.Range("I132").Value = 4147.9165
$ws.Range("K132").Value = 12443.7495
$ws.Range("M132").Value = -9913.749500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 117647624
$ws.Range("I4").Value = 166667170
$ws.Range("K4").Value = 500001510
$ws.Range("M4").Value = -500001398

$ws.Range("H12").Value = 1155.5555
$ws.Range("J12").Value = 1484.8572
$ws.Range("L12").Value = 4454.571599999999
$ws.Range("N12").Value = -4800.571599999999

$ws.Range("H13").Value = 57.875
$ws.Range("I13").Value = 24.666666
$ws.Range("J13").Value = 77.8
$ws.Range("K13").Value = 73.99999800000001
$ws.Range("L13").Value = 233.4
$ws.Range("M13").Value = 94.00000199999999
$ws.Range("N13").Value = -569.4

$ws.Range("H36").Value = 316
$ws.Range("J36").Value = 500
$ws.Range("L36").Value = 1500
$ws.Range("N36").Value = -1838

$ws.Range("H38").Value = 118
$ws.Range("I38").Value = 41
$ws.Range("J38").Value = 156.5
$ws.Range("K38").Value = 123
$ws.Range("L38").Value = 469.5
$ws.Range("M38").Value = 224
$ws.Range("N38").Value = -1163.5

$ws.Range("H97").Value = 724.5
$ws.Range("I97").Value = 399
$ws.Range("K97").Value = 1197
$ws.Range("M97").Value = -701

$ws.Range("H98").Value = 455.33334
$ws.Range("I98").Value = 299.5
$ws.Range("J98").Value = 533.25
$ws.Range("K98").Value = 898.5
$ws.Range("L98").Value = 1599.75
$ws.Range("M98").Value = 599.5
$ws.Range("N98").Value = -4595.75

$ws.Range("H103").Value = 2303.125
$ws.Range("I103").Value = 612.5
$ws.Range("J103").Value = 2866.6667
$ws.Range("K103").Value = 1837.5
$ws.Range("L103").Value = 8600.000100000001
$ws.Range("M103").Value = -958.5
$ws.Range("N103").Value = -10358.0001

$ws.Range("H122").Value = 794.7
$ws.Range("I122").Value = 331.5
$ws.Range("J122").Value = 1489.5
$ws.Range("K122").Value = 2983.5
$ws.Range("L122").Value = 13405.5
$ws.Range("M122").Value = -533.5
$ws.Range("N122").Value = -18305.5

$ws.Range("H129").Value = 1916.9166
$ws.Range("J129").Value = 2960.6
$ws.Range("L129").Value = 8881.799999999999
$ws.Range("N129").Value = -18881.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1861.3334
$ws.Range("I113").Value = 1679
$ws.Range("K113").Value = 1679
$ws.Range("M113").Value = 491

$ws.Range("H132").Value = 7999.8
$ws.Range("I132").Value = 7999.8
$ws.Range("K132").Value = 23999.4
$ws.Range("M132").Value = -21469.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4138
$ws.Range("I122").Value = 4087.6667
$ws.Range("J122").Value = 4258.8
$ws.Range("K122").Value = 12263.0001
$ws.Range("L122").Value = 12776.4
$ws.Range("M122").Value = -9813.000100000001
$ws.Range("N122").Value = -17676.4

$ws.Range("H132").Value = 2579.3125
$ws.Range("I132").Value = 2424.5454
$ws.Range("K132").Value = 7273.6362
$ws.Range("M132").Value = -4743.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 15710.143
$ws.Range("J14").Value = 14161.833
$ws.Range("L14").Value = 14161.833
$ws.Range("N14").Value = -14497.833

$ws.Range("H117").Value = 40999.5
$ws.Range("J117").Value = 40999.5
$ws.Range("L117").Value = 40999.5
$ws.Range("N117").Value = -50177.5

$ws.Range("H123").Value = 48332.332
$ws.Range("J123").Value = 48332.332
$ws.Range("L123").Value = 48332.332
$ws.Range("N123").Value = -58132.332

$ws.Range("H136").Value = 1525.6129
$ws.Range("I136").Value = 1159.579
$ws.Range("J136").Value = 2105.1667
$ws.Range("K136").Value = 3478.737
$ws.Range("L136").Value = 6315.500100000001
$ws.Range("M136").Value = -928.7370000000001
$ws.Range("N136").Value = -11415.5001
